$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.9
$ws.Range("I2").Value = 4.5
$ws.Range("J2").Value = 2.63
$ws.Range("L2").Value = 5
$ws.Range("Q2").Value = 2.4
$ws.Range("R2").Value = 1.53
$ws.Range("Y2").Value = 9
$ws.Range("AB2").Value = 34
$ws.Range("AC2").Value = 7
$ws.Range("AE2").Value = 19
$ws.Range("AF2").Value = 67
$ws.Range("AI2").Value = 21
$ws.Range("AN2").Value = 3.75
$ws.Range("AO2").Value = 11
$ws.Range("AU2").Value = 9
$ws.Range("AV2").Value = 67
$ws.Range("AX2").Value = 26
